$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.5
$ws.Range("I2").Value = 1.73
$ws.Range("J2").Value = 5.5
$ws.Range("L2").Value = 2.4
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 3
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("Z2").Value = 51
$ws.Range("AA2").Value = 41
$ws.Range("AE2").Value = 19
$ws.Range("AF2").Value = 67
$ws.Range("AK2").Value = 13
$ws.Range("AN2").Value = 6.5
$ws.Range("AO2").Value = 29
$ws.Range("AQ2").Value = 101
$ws.Range("AS2").Value = 351
$ws.Range("AU2").Value = 9
$ws.Range("AX2").Value = 9.5
$ws.Range("BC2").Value = 151
$ws.Range("BD2").Value = 151

# Row 5
$ws.Range("M5").Value = 1.13
$ws.Range("N5").Value = 6
$ws.Range("Q5").Value = 2.88
$ws.Range("R5").Value = 1.4

# Row 7
$ws.Range("U7").Value = 1.67

# Row 8
$ws.Range("G8").Value = 1.55
$ws.Range("H8").Value = 3.9
$ws.Range("I8").Value = 6.25
$ws.Range("J8").Value = 2.1
$ws.Range("L8").Value = 6
$ws.Range("O8").Value = 1.3
$ws.Range("P8").Value = 3.5
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 1.85
$ws.Range("V8").Value = 1.73
$ws.Range("W8").Value = 6
$ws.Range("Z8").Value = 11
$ws.Range("AA8").Value = 13
$ws.Range("AD8").Value = 7.5
$ws.Range("AF8").Value = 67
$ws.Range("AG8").Value = 451
$ws.Range("AJ8").Value = 19
$ws.Range("AK8").Value = 67
$ws.Range("AL8").Value = 51
$ws.Range("AN8").Value = 3.4
$ws.Range("AO8").Value = 8
$ws.Range("AQ8").Value = 26
$ws.Range("AW8").Value = 7.5
$ws.Range("AX8").Value = 34
$ws.Range("AZ8").Value = 126
$ws.Range("BB8").Value = 351
$ws.Range("BC8").Value = 151

# Row 9
$ws.Range("G9").Value = 4.33
$ws.Range("H9").Value = 3.5
$ws.Range("I9").Value = 1.83
$ws.Range("K9").Value = 2.05
$ws.Range("L9").Value = 2.5
$ws.Range("O9").Value = 1.4
$ws.Range("P9").Value = 3
$ws.Range("Q9").Value = 2.2
$ws.Range("R9").Value = 1.67
$ws.Range("U9").Value = 2
$ws.Range("V9").Value = 1.73
$ws.Range("W9").Value = 10
$ws.Range("AC9").Value = 8
$ws.Range("AD9").Value = 7
$ws.Range("AI9").Value = 7.5
$ws.Range("AP9").Value = 34
$ws.Range("AS9").Value = 301
$ws.Range("AX9").Value = 10
$ws.Range("AZ9").Value = 34
$ws.Range("BA9").Value = 51

# Row 10
$ws.Range("G10").Value = 1.85
$ws.Range("H10").Value = 3.6
$ws.Range("I10").Value = 4.1
$ws.Range("J10").Value = 2.5
$ws.Range("L10").Value = 4.33
$ws.Range("M10").Value = 1.05
$ws.Range("N10").Value = 11
$ws.Range("O10").Value = 1.25
$ws.Range("P10").Value = 4
$ws.Range("Q10").Value = 1.8
$ws.Range("R10").Value = 2
$ws.Range("U10").Value = 1.67
$ws.Range("V10").Value = 2.1
$ws.Range("W10").Value = 8.5
$ws.Range("X10").Value = 9.5
$ws.Range("Z10").Value = 17
$ws.Range("AA10").Value = 15
$ws.Range("AE10").Value = 13
$ws.Range("AG10").Value = 151
$ws.Range("AH10").Value = 13
$ws.Range("AI10").Value = 21
$ws.Range("AJ10").Value = 13
$ws.Range("AK10").Value = 41
$ws.Range("AL10").Value = 29
$ws.Range("AM10").Value = 34
$ws.Range("AN10").Value = 4
$ws.Range("AO10").Value = 10
$ws.Range("AQ10").Value = 34
$ws.Range("AR10").Value = 51
$ws.Range("AU10").Value = 7.5
$ws.Range("AW10").Value = 6
$ws.Range("AX10").Value = 21
$ws.Range("AY10").Value = 26
$ws.Range("AZ10").Value = 67
$ws.Range("BB10").Value = 151

# Row 11
$ws.Range("G11").Value = 3
$ws.Range("I11").Value = 2.63
$ws.Range("J11").Value = 3.6
$ws.Range("L11").Value = 3.4
$ws.Range("W11").Value = 7.5
$ws.Range("X11").Value = 13
$ws.Range("Y11").Value = 11
$ws.Range("Z11").Value = 29
$ws.Range("AA11").Value = 26
$ws.Range("AD11").Value = 5.5
$ws.Range("AF11").Value = 51
$ws.Range("AG11").Value = 401
$ws.Range("AH11").Value = 7
$ws.Range("AI11").Value = 12
$ws.Range("AJ11").Value = 11
$ws.Range("AK11").Value = 26
$ws.Range("AN11").Value = 4.75
$ws.Range("AO11").Value = 17
$ws.Range("AQ11").Value = 51
$ws.Range("AR11").Value = 81
$ws.Range("AS11").Value = 251
$ws.Range("AW11").Value = 4.5
